$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 108.181816
$ws.Range("I8").Value = 59
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 177
$ws.Range("L8").Value = 1800
$ws.Range("M8").Value = -38
$ws.Range("N8").Value = -2078
$ws.Range("H15").Value = 149.43
$ws.Range("I15").Value = 149.43
$ws.Range("K15").Value = 448.29
$ws.Range("M15").Value = -279.29
$ws.Range("H64").Value = 4824.5884
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 4938.625
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 4938.625
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -5434.625
$ws.Range("H67").Value = 4824.5884
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 4938.625
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 4938.625
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -6654.625
$ws.Range("H74").Value = 3575.8215
$ws.Range("I74").Value = 3100.25
$ws.Range("J74").Value = 3932.5
$ws.Range("K74").Value = 3100.25
$ws.Range("L74").Value = 3932.5
$ws.Range("M74").Value = -2164.25
$ws.Range("N74").Value = -5804.5
$ws.Range("H76").Value = 3469.3333
$ws.Range("I76").Value = 3274.7144
$ws.Range("J76").Value = 3858.5715
$ws.Range("K76").Value = 3274.7144
$ws.Range("L76").Value = 3858.5715
$ws.Range("M76").Value = -2959.7144
$ws.Range("N76").Value = -4488.5715
$ws.Range("H77").Value = 3575.8215
$ws.Range("I77").Value = 3100.25
$ws.Range("J77").Value = 3932.5
$ws.Range("K77").Value = 15501.25
$ws.Range("L77").Value = 19662.5
$ws.Range("M77").Value = -10821.25
$ws.Range("N77").Value = -29022.5
$ws.Range("H79").Value = 3469.3333
$ws.Range("I79").Value = 3274.7144
$ws.Range("J79").Value = 3858.5715
$ws.Range("K79").Value = 3274.7144
$ws.Range("L79").Value = 3858.5715
$ws.Range("M79").Value = -2182.7144
$ws.Range("N79").Value = -6042.5715
$ws.Range("H128").Value = 50400
$ws.Range("J128").Value = 50250
$ws.Range("L128").Value = 50250
$ws.Range("N128").Value = -60210
$ws.Range("H132").Value = 579078.1
$ws.Range("I132").Value = 2407.4055
$ws.Range("J132").Value = 4458499.5
$ws.Range("K132").Value = 7222.2165
$ws.Range("L132").Value = 13375498.5
$ws.Range("M132").Value = -4692.2165
$ws.Range("N132").Value = -13380558.5
$ws.Range("H138").Value = 3850410
$ws.Range("I138").Value = 3640.6667
$ws.Range("J138").Value = 5004441
$ws.Range("K138").Value = 10922.0001
$ws.Range("L138").Value = 15013323
$ws.Range("M138").Value = -5782.000100000001
$ws.Range("N138").Value = -15023603

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3749.7222
$ws.Range("I63").Value = 3746.8235
$ws.Range("K63").Value = 3746.8235
$ws.Range("M63").Value = -3060.8235
$ws.Range("H66").Value = 3749.7222
$ws.Range("I66").Value = 3746.8235
$ws.Range("K66").Value = 18734.1175
$ws.Range("M66").Value = -15302.1175
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = $null
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = $null
$ws.Range("N71").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 21741268
$ws.Range("I105").Value = 33334844
$ws.Range("J105").Value = 3311.25
$ws.Range("K105").Value = 33334844
$ws.Range("L105").Value = 3311.25
$ws.Range("M105").Value = -33333097
$ws.Range("N105").Value = -6805.25
$ws.Range("H134").Value = 8335970.5
$ws.Range("I134").Value = 2367.7144
$ws.Range("J134").Value = 27781044
$ws.Range("K134").Value = 7103.1432
$ws.Range("L134").Value = 83343132
$ws.Range("M134").Value = -4568.1432
$ws.Range("N134").Value = -83348202

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43666.113
$ws.Range("J31").Value = 82920
$ws.Range("L31").Value = 82920
$ws.Range("N31").Value = -83510
$ws.Range("H34").Value = 43666.113
$ws.Range("J34").Value = 82920
$ws.Range("L34").Value = 82920
$ws.Range("N34").Value = -83324
$ws.Range("H58").Value = 62501300
$ws.Range("I58").Value = 62501300
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 62501300
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = $null
$ws.Range("N58").Value = -62501097
$ws.Range("H99").Value = 1852.8572
$ws.Range("I99").Value = 2090
$ws.Range("J99").Value = 1675
$ws.Range("K99").Value = 2090
$ws.Range("L99").Value = 1675
$ws.Range("M99").Value = -592
$ws.Range("N99").Value = -4671
$ws.Range("H107").Value = 376.8125
$ws.Range("I107").Value = 243.6
$ws.Range("J107").Value = 598.8333
$ws.Range("K107").Value = 243.6
$ws.Range("L107").Value = 598.8333
$ws.Range("M107").Value = 1676.4
$ws.Range("N107").Value = -4438.8333
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H126").Value = 1852.8572
$ws.Range("I126").Value = 2090
$ws.Range("J126").Value = 1675
$ws.Range("K126").Value = 6270
$ws.Range("L126").Value = 5025
$ws.Range("M126").Value = -3800
$ws.Range("N126").Value = -9965
$ws.Range("H136").Value = 62501300
$ws.Range("I136").Value = 62501300
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 187503900
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = -187501350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 828.26666
$ws.Range("I34").Value = 531.8182
$ws.Range("J34").Value = 999.8946999999999
$ws.Range("K34").Value = 1595.4546
$ws.Range("L34").Value = 2999.6841
$ws.Range("M34").Value = -1511.4546
$ws.Range("N34").Value = -3167.6841
$ws.Range("H39").Value = 5966.6665
$ws.Range("J39").Value = 5966.6665
$ws.Range("L39").Value = 17899.9995
$ws.Range("N39").Value = -18487.9995
$ws.Range("H55").Value = 2000
$ws.Range("J55").Value = 2750
$ws.Range("L55").Value = 8250
$ws.Range("N55").Value = -8604
$ws.Range("H132").Value = 2172.5334
$ws.Range("I132").Value = 1448.1177
$ws.Range("K132").Value = 13033.0593
$ws.Range("M132").Value = -10503.0593

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 16000
$ws.Range("J69").Value = 16000
$ws.Range("L69").Value = 16000
$ws.Range("N69").Value = -17498
$ws.Range("H70").Value = 40607.25
$ws.Range("I70").Value = 57405.42
$ws.Range("J70").Value = 5144.4443
$ws.Range("K70").Value = 57405.42
$ws.Range("L70").Value = 5144.4443
$ws.Range("M70").Value = -57135.42
$ws.Range("N70").Value = -5684.4443
$ws.Range("H72").Value = 16000
$ws.Range("J72").Value = 16000
$ws.Range("L72").Value = 48000
$ws.Range("N72").Value = -55488
$ws.Range("H73").Value = 40607.25
$ws.Range("I73").Value = 57405.42
$ws.Range("J73").Value = 5144.4443
$ws.Range("K73").Value = 57405.42
$ws.Range("L73").Value = 5144.4443
$ws.Range("M73").Value = -56469.42
$ws.Range("N73").Value = -7016.4443
$ws.Range("H80").Value = 4469.1177
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 4468.1816
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 4468.1816
$ws.Range("M80").Value = -3502
$ws.Range("N80").Value = -6464.1816
$ws.Range("H83").Value = 4469.1177
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 4468.1816
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 22340.908
$ws.Range("M83").Value = -17508
$ws.Range("N83").Value = -32324.908
$ws.Range("H102").Value = 2692.5334
$ws.Range("I102").Value = 2779.4546
$ws.Range("J102").Value = 2453.5
$ws.Range("K102").Value = 2779.4546
$ws.Range("L102").Value = 2453.5
$ws.Range("M102").Value = -1157.4546
$ws.Range("N102").Value = -5697.5
$ws.Range("H122").Value = 2350.125
$ws.Range("I122").Value = 1855.3334
$ws.Range("J122").Value = 2986.2856
$ws.Range("K122").Value = 5566.0002
$ws.Range("L122").Value = 8958.856800000001
$ws.Range("M122").Value = -3116.0002
$ws.Range("N122").Value = -13858.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4329581.5
$ws.Range("I46").Value = 6061134.5
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 6061134.5
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = -6060946.5
$ws.Range("N46").Value = -1076
$ws.Range("H122").Value = 4320.56
$ws.Range("I122").Value = 5079.4
$ws.Range("J122").Value = 3814.6667
$ws.Range("K122").Value = 15238.2
$ws.Range("L122").Value = 11444.0001
$ws.Range("M122").Value = -12788.2
$ws.Range("N122").Value = -16344.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 29900
$ws.Range("J69").Value = 29900
$ws.Range("L69").Value = 29900
$ws.Range("N69").Value = -31398
$ws.Range("H72").Value = 29900
$ws.Range("J72").Value = 29900
$ws.Range("L72").Value = 89700
$ws.Range("N72").Value = -97188
$ws.Range("H122").Value = 2275.6191
$ws.Range("I122").Value = 1949.875
$ws.Range("J122").Value = 3318
$ws.Range("K122").Value = 5849.625
$ws.Range("L122").Value = 9954
$ws.Range("M122").Value = -3399.625
$ws.Range("N122").Value = -14854
$ws.Range("H132").Value = 48109.93
$ws.Range("I132").Value = 39604.92
$ws.Range("J132").Value = 61117.59
$ws.Range("K132").Value = 118814.76
$ws.Range("L132").Value = 183352.77
$ws.Range("M132").Value = -116284.76
$ws.Range("N132").Value = -188412.77
$ws.Range("H133").Value = 26243
$ws.Range("J133").Value = 26243
$ws.Range("L133").Value = 26243
$ws.Range("N133").Value = -36363
